$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''69.078.14'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.68%  '

$ws.Range("D3").Value = '''3.669.37'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.19%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '''676.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.30%  '

$ws.Range("D6").Value = '''160.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.32%  '

$ws.Range("D7").Value = '''3.671.18'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.14%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '''0.483'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.68%  '

$ws.Range("D10").Value = '''0.148'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.43%  '

$ws.Range("D11").Value = '''7.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.67%  '

$ws.Range("D12").Value = '''0.447'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.51%  '

$ws.Range("D13").Value = '''0.0000229'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -9.94%  '

$ws.Range("D14").Value = '''4.288.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.17%  '

$ws.Range("D15").Value = '''32.87'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -10.00%  '

$ws.Range("D16").Value = '''3.658.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.79%  '

$ws.Range("D17").Value = '''69.011.51'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.80%  '

$ws.Range("E18").Value = '  -2.11%  '

$ws.Range("D19").Value = '''16.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.90%  '

$ws.Range("D20").Value = '''6.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -9.74%  '

$ws.Range("D21").Value = '''478.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.31%  '

$ws.Range("D22").Value = '''9.80'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.95%  '

$ws.Range("D23").Value = '''0.659'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.13%  '

$ws.Range("D24").Value = '''78.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.47%  '

$ws.Range("D25").Value = '''3.805.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.44%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").Value = '''11.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.02%  '

$ws.Range("D28").Value = '''0.0000126'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -12.92%  '

$ws.Range("D29").Value = '''9.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -13.01%  '

$ws.Range("D30").Value = '''1.81'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -13.13%  '

$ws.Range("D31").Value = '''2.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -12.62%  '

$ws.Range("D32").Value = '''2.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.14%  '

$ws.Range("D33").Value = '''6.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -10.71%  '

$ws.Range("D34").Value = '''1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.02%  '

$ws.Range("D35").Value = '''26.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.43%  '

$ws.Range("E36").Value = '  -6.85%  '

$ws.Range("D37").Value = '''3.632.46'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.39%  '

$ws.Range("D38").Value = '''8.46'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.60%  '

$ws.Range("D39").Value = '''6.05'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.29%  '

$ws.Range("D40").Value = '''0.0921'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -9.86%  '

$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("D42").Value = '''2.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.55%  '

$ws.Range("D43").Value = '''1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.09%  '

$ws.Range("D44").Value = '''0.943'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -10.14%  '

$ws.Range("D45").Value = '''159.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.20%  '

$ws.Range("D46").Value = '''47.79'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.31%  '

$ws.Range("D47").Value = '''2.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -13.42%  '

$ws.Range("E48").Value = '  -4.73%  '

$ws.Range("D49").Value = '''0.000274'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -12.21%  '

$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").Value = '''382.74'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -10.79%  '

$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '''7.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.25%  '
